$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare formatting for the two brand new rows first (before row 5's
# own formatting is touched by the edits below). ---

# Row 7 will hold exactly what row 5 currently holds (APPLY / OXF-TOPMODEL-001),
# so copy row 5's current A:L formatting down onto row 7.
$ws.Range("A5:L5").Copy()
$ws.Range("A7:L7").PasteSpecial(-4122)

# Row 6 is a "fresh" data row styled like row 4's A:L formatting.
$ws.Range("A4:L4").Copy()
$ws.Range("A6:L6").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row 7: move the original row 5 content (APPLY / OXF-TOPMODEL-001) down ---
$ws.Range("A7").Value = "DEMO_TYPE2"
$ws.Range("B7").Value = "APPLY"
$ws.Range("D7").Value = "OXF-TOPMODEL-001"
$ws.Range("E7").Value = 43983
$ws.Range("G7").Value = 100
$ws.Range("H7").Value = -3500.12
$ws.Range("J7").Value = "RENT,RENT_DISCOUNT,RENT_DISCOUNT_FIXED"
$ws.Range("K7").Value = 44014
$ws.Range("L7").Value = 44073

# --- Row 6: second (25%) split item of the discount amendment ---
$ws.Range("A6").Value = "DEMO_TYPE2"
$ws.Range("B6").Value = "SIGNED"
$ws.Range("C6").Value = 44002
$ws.Range("D6").Value = "OXF-POISON-003"
$ws.Range("E6").Value = 44013
$ws.Range("G6").Value = 25
$ws.Range("J6").Value = "RENT,RENT_DISCOUNT,RENT_DISCOUNT_FIXED"
$ws.Range("K6").Value = 44105
$ws.Range("L6").Value = 44135

# --- Row 5: edited in place to become the first (50%) split item ---
$ws.Range("B5").Value = "SIGNED"
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = 44002
$ws.Range("D5").Value = "OXF-POISON-003"
$ws.Range("E5").Value = 44013
$ws.Range("G5").Value = 50
$ws.Range("F5").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("H5").ClearContents()
$ws.Range("K5").Value = 44075
$ws.Range("L5").Value = 44104

$excel.CutCopyMode = 0

$ws.Range("A9").Select() | Out-Null
